# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" between "2021-Q3" and "总计" with the
#   quarter's fund-holding detail rows.
# - Insert a new summary row at the top of "总计" for the 2022-Q1 quarter
#   (date/holding-count/holding-value), pushing the existing rows down.

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2021-Q3")
$oldTotal = $wb.Worksheets.Item("总计")

# --- format donors taken from the existing, untouched "2021-Q3" sheet ---
# style "2" (bold + thin border + centered) used for header row & index col
$headerStyleSrc = $q3.Range("B1")
# plain / default style (s=0), used so text that *looks* numeric (fund
# codes with leading zeros, decimal-looking strings, ...) keeps the exact
# "inline string" shape instead of picking up a quote-prefix style
$plainStyleSrc = $q3.Range("C2")

function Set-PlainText {
    param($cell, [string]$text)
    # leading apostrophe forces text entry even for numeric-looking values
    $cell.Value = "'" + $text
    $plainStyleSrc.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats - strips the quote-prefix flag
}

function Set-HeaderStyle {
    param($range)
    $headerStyleSrc.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# 1) Free up sheetId 3 (currently "总计") and re-create the two sheets in
#    the order/ids the target workbook uses: "2022-Q1" takes sheetId 3,
#    the re-added "总计" takes sheetId 4.
$oldTotal.Delete() | Out-Null

$q1 = $wb.Worksheets.Add($null, $q3)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# 2) Populate the new "2022-Q1" sheet -------------------------------------------------
Set-PlainText $q1.Cells.Item(1,2) "基金代码"
Set-PlainText $q1.Cells.Item(1,3) "基金名称"
Set-PlainText $q1.Cells.Item(1,4) "基金规模"
Set-PlainText $q1.Cells.Item(1,5) "股票总仓位"
Set-PlainText $q1.Cells.Item(1,6) "仓位占比"
Set-PlainText $q1.Cells.Item(1,7) "持有市值(亿元)"
Set-PlainText $q1.Cells.Item(1,8) "仓位排名"
Set-HeaderStyle $q1.Range("B1:H1")

$q1Rows = @(
    @("519983", "长信量化先锋混合A", "7.89", "94.47", "1.81", "0.1428", 9),
    @("590006", "中邮中小盘灵活配置混合", "2.39", "78.54", "2.44", "0.0583", 8),
    @("006430", "凯石澜龙头经济定期开放混合", "1.98", "69.36", "2.35", "0.0465", 9),
    @("161037", "富国中证高端制造指数增强（LOF）", "1.47", "94.32", "1.53", "0.0225", 10),
    @("004917", "中银证券祥瑞混合A", "0.10", "79.01", "2.44", "0.0024", 6),
    @("006857", "蜂巢卓睿灵活配置混合A", "0.14", "78.16", "1.22", "0.0017", 8),
    @("004918", "中银证券祥瑞混合C", "0.07", "79.01", "2.44", "0.0017", 6),
    @("004221", "长信量化先锋混合C", "0.03", "94.47", "1.81", "0.0005", 9),
    @("006858", "蜂巢卓睿灵活配置混合C", "0.04", "78.16", "1.22", "0.0005", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    Set-PlainText $q1.Cells.Item($r, 2) $row[0]
    Set-PlainText $q1.Cells.Item($r, 3) $row[1]
    Set-PlainText $q1.Cells.Item($r, 4) $row[2]
    Set-PlainText $q1.Cells.Item($r, 5) $row[3]
    Set-PlainText $q1.Cells.Item($r, 6) $row[4]
    Set-PlainText $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}
Set-HeaderStyle $q1.Range("A2:A10")

# 3) Populate the recreated "总计" sheet -----------------------------------------------
Set-PlainText $total.Cells.Item(1,2) "日期"
Set-PlainText $total.Cells.Item(1,3) "持有数量(只)"
Set-PlainText $total.Cells.Item(1,4) "持有市值(亿元)"
Set-HeaderStyle $total.Range("B1:D1")

$totalRows = @(
    @("2022-Q1", 9, 0.28),
    @("2021-Q3", 2, 0.07000000000000001),
    @("2021-Q2", 1, 0.21)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    Set-PlainText $total.Cells.Item($r, 2) $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r++
}
Set-HeaderStyle $total.Range("A2:A4")
